$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue $ws "D2" "76.607.46"
Set-TextValue $ws "E2" "  +0.95%  "
Set-TextValue $ws "D3" "3.030.09"
Set-TextValue $ws "E3" "  +4.36%  "
Set-TextValue $ws "D5" "202.03"
Set-TextValue $ws "E5" "  +1.27%  "
Set-TextValue $ws "D6" "632.53"
Set-TextValue $ws "E6" "  +6.10%  "
Set-TextValue $ws "D7" "0.999"
Set-TextValue $ws "E7" "  +0.05%  "
Set-TextValue $ws "E8" "  +1.01%  "
Set-TextValue $ws "E9" "  +6.52%  "
Set-TextValue $ws "D10" "3.027.68"
Set-TextValue $ws "E10" "  +4.41%  "
Set-TextValue $ws "D11" "0.439"
Set-TextValue $ws "E11" "  +2.80%  "
Set-TextValue $ws "E12" "  -0.36%  "
Set-TextValue $ws "D13" "5.20"
Set-TextValue $ws "E13" "  +7.38%  "
Set-TextValue $ws "D14" "3.583.76"
Set-TextValue $ws "E14" "  +4.37%  "
Set-TextValue $ws "D15" "29.61"
Set-TextValue $ws "E15" "  +7.70%  "
Set-TextValue $ws "D16" "76.528.70"
Set-TextValue $ws "E16" "  +0.96%  "
Set-TextValue $ws "D17" "0.0000195"
Set-TextValue $ws "E17" "  +2.30%  "
Set-TextValue $ws "D18" "3.023.16"
Set-TextValue $ws "E18" "  +4.20%  "
Set-TextValue $ws "D19" "13.50"
Set-TextValue $ws "E19" "  +5.32%  "
Set-TextValue $ws "D20" "8.86"
Set-TextValue $ws "E20" "  -0.45%  "
Set-TextValue $ws "D21" "378.12"
Set-TextValue $ws "E21" "  +1.97%  "
Set-TextValue $ws "D22" "2.32"
Set-TextValue $ws "E22" "  +0.80%  "
Set-TextValue $ws "D23" "4.39"
Set-TextValue $ws "E23" "  +3.06%  "
Set-TextValue $ws "D24" "73.77"
Set-TextValue $ws "E24" "  +4.07%  "
Set-TextValue $ws "D25" "3.184.23"
Set-TextValue $ws "D26" "4.41"
Set-TextValue $ws "E26" "  +5.94%  "
Set-TextValue $ws "E27" "  +0.07%  "
Set-TextValue $ws "D28" "10.01"
Set-TextValue $ws "E28" "  +4.10%  "
Set-TextValue $ws "D29" "0.0000112"
Set-TextValue $ws "E29" "  +3.93%  "
Set-TextValue $ws "D30" "1.00"
Set-TextValue $ws "E30" "  +0.00%  "
Set-TextValue $ws "D31" "8.37"
Set-TextValue $ws "E31" "  +8.73%  "
Set-TextValue $ws "D32" "1.43"
Set-TextValue $ws "E32" "  +1.94%  "
Set-TextValue $ws "D33" "517.64"
Set-TextValue $ws "E33" "  +3.35%  "
Set-TextValue $ws "D34" "1.97"
Set-TextValue $ws "E34" "  +9.17%  "
Set-TextValue $ws "D35" "0.999"
Set-TextValue $ws "E35" "  +0.08%  "
Set-TextValue $ws "D36" "20.82"
Set-TextValue $ws "E36" "  +3.81%  "
Set-TextValue $ws "D37" "163.85"
Set-TextValue $ws "E37" "  -1.04%  "
Set-TextValue $ws "D38" "0.385"
Set-TextValue $ws "E38" "  +11.98%  "
Set-TextValue $ws "E39" "  +1.96%  "
Set-TextValue $ws "E40" "  +6.91%  "
Set-TextValue $ws "D41" "187.75"
Set-TextValue $ws "E41" "  +4.27%  "
Set-TextValue $ws "E42" "  +0.44%  "
Set-TextValue $ws "E43" "  +0.29%  "
Set-TextValue $ws "D44" "5.19"
Set-TextValue $ws "E44" "  +5.02%  "
Set-TextValue $ws "D45" "42.38"
Set-TextValue $ws "E45" "  +5.66%  "
Set-TextValue $ws "E46" "  +7.54%  "
Set-TextValue $ws "D47" "1.69"
Set-TextValue $ws "E47" "  +2.95%  "
Set-TextValue $ws "D48" "2.49"
Set-TextValue $ws "E48" "  +6.97%  "
Set-TextValue $ws "D49" "0.723"
Set-TextValue $ws "E49" "  +10.66%  "
Set-TextValue $ws "E50" "  +6.98%  "
Set-TextValue $ws "E51" "  +5.80%  "
